# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain-text values (e.g. "43.709.91", "  +1.23%  ")
# that Excel would otherwise auto-convert to numbers. Force them to text first,
# then restore the original (default) cell style so no stray formatting is left behind.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '43.709.91'
$ws.Range("D3").Value = '2.283.35'
$ws.Range("E3").Value = '  +1.23%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '119.99'
$ws.Range("E5").Value = '  +7.65%  '
$ws.Range("D6").Value = '267.65'
$ws.Range("E6").Value = '  +1.85%  '
$ws.Range("D7").Value = '0.648'
$ws.Range("E7").Value = '  +5.35%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  +5.91%  '
$ws.Range("D10").Value = '48.45'
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '9.27'
$ws.Range("E12").Value = '  +6.90%  '
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").Value = '15.75'
$ws.Range("E14").Value = '  +2.87%  '
$ws.Range("D15").Value = '0.922'
$ws.Range("E15").Value = '  +8.66%  '
$ws.Range("D16").Value = '2.625.90'
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D17").Value = '2.284.33'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("D18").Value = '43.782.05'
$ws.Range("E18").Value = '  +1.77%  '
$ws.Range("D19").Value = '0.0000111'
$ws.Range("E19").Value = '  +4.13%  '
$ws.Range("D20").Value = '6.97'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").Value = '72.44'
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("E22").Value = '  +1.62%  '
$ws.Range("D23").Value = '236.44'
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("D24").Value = '9.59'
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").Value = '2.89'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("D26").Value = '12.08'
$ws.Range("E26").Value = '  +7.77%  '
$ws.Range("D28").Value = '43.15'
$ws.Range("E28").Value = '  +4.69%  '
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").Value = '172.84'
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("D32").Value = '21.81'
$ws.Range("E32").Value = '  +3.06%  '
$ws.Range("D33").Value = '0.0929'
$ws.Range("E33").Value = '  +4.11%  '
$ws.Range("E34").Value = '  +4.89%  '
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").Value = '  +5.03%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.28'
$ws.Range("E36").Value = '  +14.91%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.0386'
$ws.Range("E37").Value = '  +11.00%  '
$ws.Range("D38").Value = '4.63'
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = '0.109'
$ws.Range("E39").Value = '  +5.68%  '
$ws.Range("D40").Value = '2.58'
$ws.Range("E40").Value = '  +8.24%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").Value = '13.90'
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '74.35'
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").Value = '0.240'
$ws.Range("E43").Value = '  +3.33%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '1.40'
$ws.Range("E45").Value = '  +2.70%  '
$ws.Range("D46").Value = '5.87'
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("D47").Value = '74.66'
$ws.Range("E47").Value = '  +44.46%  '
$ws.Range("D48").Value = '1.28'
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.101'
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '8.54'
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("D51").Value = '102.45'
$ws.Range("E51").Value = '  +3.14%  '

$priceVolRange.Style = "Normal"

